# Daily attendance processing - 2025-11-30 10:26:29
# Normalize the "Recorded By" column (G) so that the literal token
# "System" (case-sensitive, exact match) is always listed first among
# the comma-separated recorder names, preserving the relative order of
# the remaining entries (e.g. a lowercase "system" stays in place, just
# after the canonical "System" entry is moved to the front).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ",\s*"
        $rest = @()
        $foundSystem = $false

        foreach ($part in $parts) {
            if ($part.Equals("System")) {
                $foundSystem = $true
            } else {
                $rest += $part
            }
        }

        if ($foundSystem) {
            $ordered = @("System") + $rest
            $newValue = $ordered -join ", "
            $cell.Value2 = $newValue
        }
    }
}
